$d = $word.ActiveDocument
$table = $d.Tables(1)

# Row 2: "Matching Effect" point estimates
$table.Cell(2, 2).Range.Text = "-1114.80***"
$table.Cell(2, 3).Range.Text = "-967.05***"

# Row 3: standard errors for Matching Effect
$table.Cell(3, 2).Range.Text = "(38.50)"
$table.Cell(3, 3).Range.Text = "(35.11)"

# Row 4: "(Intercept)" point estimates
$table.Cell(4, 2).Range.Text = "1883.33***"
$table.Cell(4, 3).Range.Text = "1771.27***"

# Row 5: standard errors for (Intercept)
$table.Cell(5, 2).Range.Text = "(30.54)"
$table.Cell(5, 3).Range.Text = "(27.94)"

# Row 6: Num.Obs. -- both cells currently read "5928" but must diverge
$table.Cell(6, 2).Range.Text = "5186"
$table.Cell(6, 3).Range.Text = "4888"
